# Generate Report for Handoff
# Update the localization status report: mark the file as ready for
# handoff (was "In Translation") and refresh the handoff timestamps on
# all three sheets (Overview + per-locale detail sheets).

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-08-29 22:41:24"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-08-29 22:41:19"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-08-29 22:41:24"

$wsOverview.Columns("E:F").EntireColumn.AutoFit()
$wsZhCn.Columns("C:C").EntireColumn.AutoFit()
$wsDeDe.Columns("C:C").EntireColumn.AutoFit()
